$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.485.61"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +0.25%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.797.96"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -0.23%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "317.02"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.25%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.04%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5426"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -1.05%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3783"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -1.34%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07505"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -0.72%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.83"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -2.27%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.108"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -1.36%  "

$ws.Range("E12").Value = "  +0.08%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.70"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -2.03%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.161"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -0.43%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.299"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.19%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.799.38"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -0.14%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "89.65"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -1.21%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001067"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.00%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06502"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.76%  "

$ws.Range("E20").Value = "  +1.37%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.001"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.04%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.957"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.25%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.483.34"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.29%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.11"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -1.07%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.078"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -1.68%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "159.60"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +1.14%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.46"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.70%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.000.02"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -0.51%  "

$ws.Range("E29").Value = "  -4.28%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "122.80"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -0.66%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.110"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -4.46%  "

$ws.Range("E32").Value = "  +3.08%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.615"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -1.54%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.652"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.56%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.2279"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.25%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06464"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +3.29%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02301"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.82%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.619"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -3.21%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.032"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +0.62%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6207"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -2.48%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "11.20"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -2.99%  "

$ws.Range("E42").Value = "  +4.98%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.194"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +3.38%  "

$ws.Range("E44").Value = "  -0.05%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.30"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.37%  "

$ws.Range("E46").Value = "  +0.38%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5829"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -2.18%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "127.29"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +3.04%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.205"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +5.46%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.950"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.24%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06882"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.16%  "

